$d = $word.ActiveDocument

$p1 = $d.Paragraphs(1)

# Add a paragraph border (top/left/bottom/right, space=5 from the text) to the first paragraph.
$b = $p1.Borders
$b.DistanceFromTop = 5
$b.DistanceFromLeft = 5
$b.DistanceFromBottom = 5
$b.DistanceFromRight = 5

# Increase the paragraph's left indent from 120 twips (6pt) to 225 twips (11.25pt).
$p1.LeftIndent = 11.25

# Remove the trailing stand-alone space run at the end of the first paragraph
# (the last character before the paragraph mark).
$pr = $p1.Range
$spaceRange = $d.Range($pr.End - 2, $pr.End - 1)
if ($spaceRange.Text -eq " ") {
    $spaceRange.Delete()
}

# Update the merge-field id text to reflect the new topic naming.
$d.Content.Find.Execute("**ID__AFFARS_pgi_5349_topic_4__ID**", $true, $false, $false, $false, $false,
                         $true, 1, $false, "**ID__AFFARS_SMC_PGI_5349__ID**", 2)
